# Lab 2 deflections: sort the data table (A2:D57) by the Manager_narrative
# column (B), ascending, keeping the header row (row 1) fixed, then leave
# the selection on C55 as the author did after sorting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A2:D57")
$sortKey   = $ws.Range("B2:B57")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($sortKey)
$ws.Sort.SetRange($dataRange)
$ws.Sort.Header = 2
$ws.Sort.Apply()

$ws.Range("C55").Select() | Out-Null
